# Fix errors in template input files:
#  - E4 ("file to define each conductor") should point to the template
#    conductor definition file, not the generic one.
#  - E5 ("file to define the external environment") should point to the
#    template environment input file, not the generic one.
#  - Update the active cell selection on the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TRANSIENT")

$ws.Range("E4").Value = "template_conductor_definition.xlsx"
$ws.Range("E5").Value = "template_environment_input.xlsx"

$ws.Activate()
$ws.Range("E9").Select()
